$d = $word.ActiveDocument

# wdReplaceOne = 1 so that only the first placeholder occurrence (under
# "Programador #4") is updated; the later placeholders for programmers
# #5-#10 must stay untouched.
$d.Content.Find.Execute(
    "(Escribir nombre completo y número de carnet)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Daniel Alexander Pérez Ramos - 2017145", 1)

$d.Content.Find.Execute(
    "(Escribir actividad asignada)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Creación de 3 vistas (Empleado, producto y detalle venta).", 1)
